$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Implementation Facts/Decisions" / "/Assumptions" had been
# split across two runs; collapse them back into a single run holding
# the full heading text "Implementation Facts/Decisions/Assumptions".
# A no-op Find/Replace over the concatenated text naturally merges the
# runs that make it up into one.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Implementation Facts/Decisions/Assumptions", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Implementation Facts/Decisions/Assumptions", 2)

# ---------------------------------------------------------------------
# Change 2: the red "TODO - Specify having built cpputest..." note.
# Before:
#   ...cpputest, defining the CPP_UNIT_HOME  env var, adding this to
#   linker props $(CPP_UNIT_HOME)\lib, adding $(CPP_UNIT_HOME)\include
#   to compiler props, having built the cppunit_dll project ... directory.
#   [_GoBack bookmark sits at the very end of the paragraph]
#
# After:
#   ...cpputest to the correct x86/x64[_GoBack], defining the
#   CPP_UNIT_HOME  env var, adding this to linker props
#   $(CPP_UNIT_HOME)\lib, adding $(CPP_UNIT_HOME)\include to compiler
#   props, having built the cppunit_dll project ... directory.
#   [_GoBack bookmark now sits right after "x86/x64", before ", defining"]
# ---------------------------------------------------------------------

# Remember where the _GoBack bookmark used to live, then remove it - it
# gets re-anchored further up the sentence below.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Replace the chunk that used to follow "cpputest" (", defining the ...
# compiler props") with the new "to the correct x86/x64" wording.
# (single-quoted so the literal "$(...)" text needs no escaping)
$oldChunk = ', defining the CPP_UNIT_HOME  env var, adding this to linker props $(CPP_UNIT_HOME)\lib, adding $(CPP_UNIT_HOME)\include to compiler props'
$rng = $d.Content
$rng.Find.Execute($oldChunk, $false, $false, $false, $false, $false, $true, 1, $false, " to the correct x86/x64", 2)

# $rng now spans the freshly-inserted replacement text; collapse to its
# end so we're sitting right after "x86/x64".
$rng.Collapse(0)
$anchor = $rng.Start

# Re-insert the text that used to follow immediately (", defining the
# CPP_UNIT_HOME ... compiler props") right after that point, so the
# sentence still reads correctly.
$rng.InsertAfter($oldChunk)

# Finally drop the _GoBack bookmark back in, collapsed at the original
# anchor point - i.e. right after "x86/x64" and before the
# re-inserted ", defining the ..." text.
$d.Bookmarks.Add("_GoBack", $d.Range($anchor, $anchor))

Write-Output "done"
